$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "Most Active Briefing(s):" (first occurrence, the "external"/
#    overall briefing count paragraph) gets " Externally" inserted
#    (italic) between "(s)" and ":" -- mirroring the "Internally"
#    paragraph further down in the document.
# ------------------------------------------------------------------
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("Briefing(s):", $false, $false, $false, $false, $false, `
                                  $false, 1, $false, "", 0)

if ($found) {
    # Position right after "(s)" and before ":" -- that's End-1 of the
    # found range (the ":" is the very last character matched).
    $insertPos = $findRange.End - 1

    # Insert the separating space first.
    $spaceRange = $d.Range($insertPos, $insertPos)
    $spaceRange.InsertAfter(" ")

    # Insert "Externally" right after the space, then italicise it.
    $extStart = $insertPos + 1
    $extRange = $d.Range($extStart, $extStart)
    $extRange.InsertAfter("Externally")

    $extFormatRange = $d.Range($extStart, $extStart + 10)
    $extFormatRange.Font.Italic = $true
}

# ------------------------------------------------------------------
# 2. The four "BRIEFCOUNT2..5" list paragraphs change their
#    indentation from left=2160 to left=2880 + firstLine=720.
# ------------------------------------------------------------------
$targets = @("{{BRIEFCOUNT2}}", "{{BRIEFCOUNT3}}", "{{BRIEFCOUNT4}}", "{{BRIEFCOUNT5}}")

foreach ($target in $targets) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        if ($para.Range.Text -like "*$target*") {
            $para.LeftIndent = 144
            $para.FirstLineIndent = 36
            break
        }
    }
}
